$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns H and I
$ws.Range("H1").Value = "Tempo Heuristica"
$ws.Range("I1").Value = "Tempo Total"

# Row 2 updates
$ws.Range("G2").Value = 0.05321311950683594
$ws.Range("H2").Value = 0.01083183288574219
$ws.Range("I2").Value = 0.06404495239257812

# Row 3 updates
$ws.Range("G3").Value = 0.04790878295898438
$ws.Range("H3").Value = 0.01452851295471191
$ws.Range("I3").Value = 0.06243729591369629

# Row 4 updates
$ws.Range("G4").Value = 0.04992318153381348
$ws.Range("H4").Value = 0.01286220550537109
$ws.Range("I4").Value = 0.06278538703918457
